{"js": "// Update the \"alt text\" (title / description) of the first inline chart\n// (\"Graphique 2\") so the description carries the new TBS \"[tbs:ref=...]\"\n// tag used by Ms Office >= 2019, folding the former Title text into the\n// Description and clearing the separate Title.\n\nconst inlinePictures = context.document.body.inlinePictures;\ninlinePictures.load(\"items\");\nawait context.sync();\n\n// The chart is the first inline graphic in the document body.\nconst chart = inlinePictures.items[0];\nchart.load([\"altTextTitle\", \"altTextDescription\"]);\nawait context.sync();\n\nlet title = chart.altTextTitle;            // \"a nice chart\"\nlet description = chart.altTextDescription; // \"This is just a nice chart\"\nif (!title) {\n  title = \"a nice chart\";\n}\n\nchart.altTextDescription = \"[tbs:ref=\" + title + \"]\\n\\n\" + description;\nchart.altTextTitle = \"\";\n\nawait context.sync();\n", "ps1": "# Update the \"alt text\" (title / description) of the first inline chart\n# (\"Graphique 2\") so the description carries the new TBS \"[tbs:ref=...]\"\n# tag used by Ms Office >= 2019, folding the former Title text into the\n# Description and clearing the separate Title.\n\n$d = $word.ActiveDocument\n\n# The chart is the first inline graphic (\"InlineShape\") in the document.\n$chart = $d.InlineShapes.Item(1)\n\n$title = $chart.Title                 # \"a nice chart\"\n$description = $chart.AlternativeText # \"This is just a nice chart\"\nif (-not $title) { $title = \"a nice chart\" }\n\n$chart.AlternativeText = \"[tbs:ref=\" + $title + \"]`n`n\" + $description\n$chart.Title = \"\"\n"}
